# Actualización automática 2025-07-04 14:05:09
#
# Updates sales figures for "CASTRO ALCIVAR EDA MARIA" across the three
# sheets of the workbook:
#   - "VENTAS POR GRUPO"     : sales by product group per client
#   - "VENTA MENSUAL"        : sales by month per client
#   - "CUMPLIMIENTO MENSUAL" : budget-compliance roll-up by product group

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("M29").Value = 514.08
$wsGrupo.Range("D45").Value = 89.76000000000001
$wsGrupo.Range("H45").Value = 213.3
$wsGrupo.Range("I45").Value = 103.5
$wsGrupo.Range("L45").Value = 73.91
$wsGrupo.Range("M50").Value = 921.4299999999999

# Row 55 "X de 53" compliance counters
$wsGrupo.Range("D55").Value = "1 de 53"
$wsGrupo.Range("H55").Value = "1 de 53"
$wsGrupo.Range("I55").Value = "2 de 53"
$wsGrupo.Range("L55").Value = "1 de 53"
$wsGrupo.Range("M55").Value = "4 de 53"

# ---------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F29").Value = 514.08
$wsMensual.Range("F45").Value = 492.84
$wsMensual.Range("F50").Value = 1247.67
$wsMensual.Range("F55").Value = 6379.46

# ---------------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL
# ---------------------------------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3 - 240X80 PORCELANATO
$wsCumplimiento.Range("D3").Value = 89.76000000000001
$wsCumplimiento.Range("E3").Value = 13638.24
$wsCumplimiento.Range("F3").Value = 0.006538461538461539

# Row 7 - INODOROS
$wsCumplimiento.Range("D7").Value = 213.3
$wsCumplimiento.Range("E7").Value = 2986.7
$wsCumplimiento.Range("F7").Value = 0.06665625

# Row 8 - LAVABOS
$wsCumplimiento.Range("D8").Value = 129.6
$wsCumplimiento.Range("E8").Value = 870.4
$wsCumplimiento.Range("F8").Value = 0.1296

# Row 15 - PIEDRA SINTERIZADA
$wsCumplimiento.Range("D15").Value = 73.91
$wsCumplimiento.Range("E15").Value = 20616.09
$wsCumplimiento.Range("F15").Value = 0.003572257129047849

# Row 16 - PORCELANATO
$wsCumplimiento.Range("D16").Value = 2476.28
$wsCumplimiento.Range("E16").Value = 52244.95
$wsCumplimiento.Range("F16").Value = 0.04525263777879262

# Row 19 - TOTAL
$wsCumplimiento.Range("D19").Value = 6787.289999999999
$wsCumplimiento.Range("E19").Value = 98425.57999999999
$wsCumplimiento.Range("F19").Value = 0.06451007372006864

# Column F ("CUMPLIMIENTO") widened from 25 to 26 characters.
# ColumnWidth is specified in characters of the Normal-style font and
# Excel stores it in the XML as a slightly different (pixel-rounded)
# number; 25.17 is the calibrated input that round-trips to a stored
# width of exactly 26.
$wsCumplimiento.Columns.Item(6).ColumnWidth = 25.17
